$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.657.41'
$ws.Range('E2').Value = '  +0.41%  '

$ws.Range('D3').Value = '1.597.12'
$ws.Range('E3').Value = '  +0.91%  '

$ws.Range('E4').Value = '  +0.03%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '210.73'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.01%  '

$ws.Range('E6').Value = '  +1.44%  '

$ws.Range('E7').Value = '  +0.05%  '

$ws.Range('E8').Value = '  -0.05%  '

$ws.Range('E9').Value = '  -1.00%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '19.61'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +0.59%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0843'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +1.41%  '

$ws.Range('D12').Value = '1.821.13'
$ws.Range('E12').Value = '  +0.90%  '

$ws.Range('D13').Value = '1.568.38'
$ws.Range('E13').Value = '  -1.04%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.02'
$ws.Range('D14').Style = 'Normal'

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.520'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -1.52%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '64.73'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.62%  '

$ws.Range('D17').Value = '26.642.83'
$ws.Range('E17').Value = '  +0.19%  '

$ws.Range('D18').Value = '0.0₃0728'
$ws.Range('E18').Value = '  -0.20%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '208.44'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.00%  '

$ws.Range('E20').Value = '  +0.08%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.76'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.84%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.25'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.10%  '

$ws.Range('E23').Value = '  -3.11%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '8.90'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.25%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '145.44'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.67%  '

$ws.Range('E26').Value = '  +0.07%  '

$ws.Range('E27').Value = '  -2.88%  '

$ws.Range('E28').Value = '  +2.15%  '

$ws.Range('E29').Value = '  -0.17%  '

$ws.Range('E30').Value = '  +0.78%  '

$ws.Range('E31').Value = '  +0.02%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.24'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -0.77%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.656'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -2.10%  '

$ws.Range('E34').Value = '  -0.17%  '

$ws.Range('D35').Value = '1.292.11'
$ws.Range('E35').Value = '  -2.11%  '

$ws.Range('E36').Value = '  +0.32%  '

$ws.Range('E37').Value = '  -1.03%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.0172'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -0.79%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.845'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +2.91%  '

$ws.Range('E40').Value = '  +0.00%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '5.39'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +1.66%  '

$ws.Range('E42').Value = '  +1.23%  '

$ws.Range('E43').Value = '  +0.23%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '63.69'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +1.08%  '

$ws.Range('D45').Value = '1.734.06'
$ws.Range('E45').Value = '  +0.88%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.897'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +8.00%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '90.02'
$ws.Range('D47').Style = 'Normal'

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.61'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.66%  '

$ws.Range('E49').Value = '  +1.82%  '

$ws.Range('E50').Value = '  -0.26%  '

$ws.Range('B51').Value = 'EnergySwap'
$ws.Range('C51').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '7.45'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.15%  '

